$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("sample")

# The "sample" import template had a leftover, duplicate pair of columns
# (AE: "location" / AF: "flash point") that repeated - with shorter/partial
# values - what columns AB ("location") and AC ("flash point") already
# documented. Clear them out (header + all 6 example rows) so the template
# only shows the one, fully-populated set of example columns.
$ws.Range("AE1:AF7").ClearContents()

# Now that the duplicate columns are gone, widen the real "location" column
# so the longer example values ("room x- shelf 1" ...) are fully visible
# instead of relying on the old best-fit width.
$ws.Columns("AB").ColumnWidth = 22.3

# Tidy up the leftover duplicate "Normal" cell style that the template
# carried around unused (Excel keeps only the built-in one).
$wb.Styles.Item("Normal").Delete()

# Leave the selection where it was left after reviewing the cleaned-up sheet.
$ws.Range("AC20").Select()
